$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: "time_taken" metadata (time taken to process each panel row)
# Header F1 mirrors the formatting of the other header cells (e.g. E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:41:50.764040",
    "2021-10-05 13:41:50.764053",
    "2021-10-05 13:41:50.764057",
    "2021-10-05 13:41:50.764060",
    "2021-10-05 13:41:50.764064",
    "2021-10-05 13:41:50.764067",
    "2021-10-05 13:41:50.764070",
    "2021-10-05 13:41:50.764072",
    "2021-10-05 13:41:50.764076",
    "2021-10-05 13:41:50.764079",
    "2021-10-05 13:41:50.764082",
    "2021-10-05 13:41:50.764085",
    "2021-10-05 13:41:50.764088",
    "2021-10-05 13:41:50.764090",
    "2021-10-05 13:41:50.764093",
    "2021-10-05 13:41:50.764096",
    "2021-10-05 13:41:50.764099",
    "2021-10-05 13:41:50.764103",
    "2021-10-05 13:41:50.764106",
    "2021-10-05 13:41:50.764109",
    "2021-10-05 13:41:50.764111",
    "2021-10-05 13:41:50.764114",
    "2021-10-05 13:41:50.764117",
    "2021-10-05 13:41:50.764120",
    "2021-10-05 13:41:50.764124",
    "2021-10-05 13:41:50.764127",
    "2021-10-05 13:41:50.764130",
    "2021-10-05 13:41:50.764133",
    "2021-10-05 13:41:50.764135",
    "2021-10-05 13:41:50.764138",
    "2021-10-05 13:41:50.764141",
    "2021-10-05 13:41:50.764144",
    "2021-10-05 13:41:50.764148",
    "2021-10-05 13:41:50.764151",
    "2021-10-05 13:41:50.764154",
    "2021-10-05 13:41:50.764157",
    "2021-10-05 13:41:50.764160",
    "2021-10-05 13:41:50.764163",
    "2021-10-05 13:41:50.764166",
    "2021-10-05 13:41:50.764169",
    "2021-10-05 13:41:50.764172",
    "2021-10-05 13:41:50.764175",
    "2021-10-05 13:41:50.764178",
    "2021-10-05 13:41:50.764181",
    "2021-10-05 13:41:50.764184",
    "2021-10-05 13:41:50.764187",
    "2021-10-05 13:41:50.764190",
    "2021-10-05 13:41:50.764193",
    "2021-10-05 13:41:50.764196",
    "2021-10-05 13:41:50.764199",
    "2021-10-05 13:41:50.764202",
    "2021-10-05 13:41:50.764205",
    "2021-10-05 13:41:50.764208",
    "2021-10-05 13:41:50.764211",
    "2021-10-05 13:41:50.764214",
    "2021-10-05 13:41:50.764217",
    "2021-10-05 13:41:50.764220",
    "2021-10-05 13:41:50.764223",
    "2021-10-05 13:41:50.764226",
    "2021-10-05 13:41:50.764229",
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

$excel.CutCopyMode = 0

